# Add 4 new weekly scoreboard rows (73-76) and widen column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: Participant, Date(serial), Workout Type, Total Duration,
# Total Distance, Total Elevation, Zone1..Zone5, Workout Level, Week
$rows = @(
    @{ Row=73; A="Steven"; B=45462; C="Walk"; D=29;  E=1.22;  F=30;   G=29; H=0;  I=0;  J=0;  K=0;  L="Agile Antelope"; M=2 },
    @{ Row=74; A="Matt";   B=45463; C="Ride"; D=84;  E=22.63; F=1283; G=11; H=68; I=3;  J=0;  K=0;  L="Agile Antelope"; M=2 },
    @{ Row=75; A="Steven"; B=45463; C="Run";  D=26;  E=2.32;  F=69;   G=1;  H=3;  I=10; J=11; K=0;  L="Agile Antelope"; M=2 },
    @{ Row=76; A="Steven"; B=45463; C="Walk"; D=19;  E=0.77;  F=23;   G=13; H=5;  I=2;  J=0;  K=0;  L="Agile Antelope"; M=2 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Copy the date cell's number format from an existing dated cell (B2)
    # so the new date cell reuses the same style (short-date numFmt) rather
    # than minting a new cell style.
    $ws.Cells.Item(2, 2).Copy($ws.Cells.Item($row, 2))

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# Widen column C (Workout Type) to fit the new longer entries; drop the
# previous best-fit autosize in favor of an explicit width.
$ws.Columns.Item(3).ColumnWidth = 13.5

# Move the active selection to the new last cell, matching where Excel
# would leave the cursor after entering this data.
$ws.Range("M76").Select() | Out-Null
